$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1) — reorder / expand to A1:AA1
# ---------------------------------------------------------------------------
$headers = [ordered]@{
    "A1" = "identificacion";
    "B1" = "nombre";
    "C1" = "email";
    "D1" = "estado";
    "E1" = "tipo de documento";
    "F1" = "ciudad expedición";
    "G1" = "fecha de nacimiento";
    "H1" = "dirección";
    "I1" = "teléfono";
    "J1" = "cargo";
    "K1" = "rh";
    "L1" = "método de pago";
    "M1" = "banco";
    "N1" = "tipo de cuenta";
    "O1" = "número de cuenta";
    "P1" = "salario";
    "Q1" = "base";
    "R1" = "fecha de inicio";
    "S1" = "subsidio";
    "T1" = "contratación";
    "U1" = "tipo de contrato";
    "V1" = "eps";
    "W1" = "cesantias";
    "X1" = "pensiones";
    "Y1" = "area";
    "Z1" = "caja de compensación";
    "AA1" = "arl";
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# ---------------------------------------------------------------------------
# 2. New employee record (row 2)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 1076899023
$ws.Range("B2").Value = "Armando Cristancho"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "CC"
$ws.Range("F2").Value = "soacha"
$ws.Range("G2").Value = 20
$ws.Range("H2").Value = "Carrera demo 123"
$ws.Range("I2").Value = 3098997665
$ws.Range("J2").Value = "CONDUCTOR"
$ws.Range("K2").Value = "O+"
$ws.Range("L2").Value = "EFECTIVO"
$ws.Range("M2").Value = "BANCOLOMBIA"
$ws.Range("N2").Value = "AHORROS"
$ws.Range("O2").Value = 98888288828
$ws.Range("P2").Value = 8000000
$ws.Range("Q2").Value = "MENSUAL"
$ws.Range("S2").Value = 100000
$ws.Range("T2").Value = "EMPLEADO"
$ws.Range("U2").Value = "INDEFINIDO"
$ws.Range("V2").Value = "NUEVA EPS"
$ws.Range("W2").Value = "PORVENIR"
$ws.Range("X2").Value = "PORVENIR"
$ws.Range("Y2").Value = "OPERATIVO"
$ws.Range("Z2").Value = "COMFACASANARE"
$ws.Range("AA2").Value = "POSITIVA COMPAÑÍA DE SEGUROS"

# Email cell becomes a mailto hyperlink (also applies the built-in Hyperlink style)
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:actancho@dem.com", [Type]::Missing, [Type]::Missing, "actancho@dem.com") | Out-Null

# Date of hire (R2) — set the number format first so no transient custom
# numFmt entry is created, then write a pure date (no time-of-day part).
$ws.Range("R2").NumberFormat = "mm-dd-yy"
$ws.Range("R2").Value = (Get-Date -Year 2021 -Month 2 -Day 25 -Hour 0 -Minute 0 -Second 0)

# ---------------------------------------------------------------------------
# 3. Column widths (best-fit widths for the populated columns)
# ---------------------------------------------------------------------------
$widths = [ordered]@{
    1 = 13.140625;
    2 = 8;
    3 = 6;
    4 = 7;
    5 = 18;
    6 = 17.28515625;
    7 = 19.140625;
    8 = 9.140625;
    9 = 11;
    10 = 5.7109375;
    11 = 2.85546875;
    12 = 15.42578125;
    13 = 6.28515625;
    14 = 13.7109375;
    15 = 17.28515625;
    16 = 8;
    17 = 5.140625;
    18 = 13.85546875;
    19 = 8.42578125;
    20 = 12;
    21 = 15.140625;
    22 = 4.140625;
    23 = 9.140625;
    24 = 10.140625;
    25 = 4.85546875;
    26 = 20.42578125;
    27 = 3.28515625;
}
foreach ($col in $widths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $widths[$col]
}

# ---------------------------------------------------------------------------
# 4. Selection
# ---------------------------------------------------------------------------
$ws.Range("Z3").Select() | Out-Null
